# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (between "2021-Q4" and "总计") holding
#   the quarter's fund-holding detail rows, in the same layout as the
#   "2021-Q4" sheet.
# - Prepend a 2022-Q1 summary row to the "总计" sheet, pushing the existing
#   2021-Q4 summary row down one row.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" worksheet right after "2021-Q4" -----------
# (worksheet references captured before a Worksheets.Add() can go stale, so
# resolve "2021-Q4" fresh right before using it as the anchor.)
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("2021-Q4"))
$newSheet.Name = "2022-Q1"

# Re-resolve "2021-Q4" again (post sheet-insert) as the format donor.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# Clone header formatting (bold/centered/bordered) + the index-column style.
$q4Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q4Sheet.Range("A2").Copy($newSheet.Range("A2"))

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row (kept as text to match the source formatting, except the rank
# which is numeric). Force text entry via NumberFormat, then drop the
# number format again so the cell is left with no explicit style, matching
# the unstyled data rows elsewhere in the workbook.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "007835"
$newSheet.Range("C2").Value = "国泰鑫睿混合"
$newSheet.Range("D2").Value = "9.37"
$newSheet.Range("E2").Value = "78.94"
$newSheet.Range("F2").Value = "3.09"
$newSheet.Range("G2").Value = "0.2895"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# --- 2. Prepend a 2022-Q1 row to the "总计" summary sheet -------------------
$sumSheet = $wb.Worksheets.Item("总计")
$sumSheet.Range("A2:D2").Insert()

# The old row 2 (2021-Q4) is now row 3; clone its index-column style for the
# new row 2, and drop the formatting the insert auto-extended onto B2:D2 so
# the new data row is unstyled, like the rest of the sheet's data rows.
$sumSheet.Range("A3").Copy($sumSheet.Range("A2"))
$sumSheet.Range("B2:D2").ClearFormats()

$sumSheet.Range("A2").Value = 0
$sumSheet.Range("B2").Value = "2022-Q1"
$sumSheet.Range("C2").Value = 1
$sumSheet.Range("D2").Value = 0.29

$sumSheet.Range("A3").Value = 1

# Restore the originally-active sheet/selection.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Activate() | Out-Null
$q4Sheet.Range("A1").Select() | Out-Null
